$d = $word.ActiveDocument

$d.Content.Find.Execute("[Music]", $true, $false, $false, $false, $false, $true, 1, $false, "[ሙዚቃ]", 2) | Out-Null
$d.Content.Find.Execute("them Fil and Mike who meet each other", $true, $false, $false, $false, $false, $true, 1, $false, "ፊል እና ሚኪ ተገናኙ ", 2) | Out-Null
$d.Content.Find.Execute("again after a long time. After some", $true, $false, $false, $false, $false, $true, 1, $false, "ከረጅም ጊዜ በኋላ፡፡ ከተወሰነ", 2) | Out-Null
$d.Content.Find.Execute("chatting, Phil says he has three children, then", $true, $false, $false, $false, $false, $true, 1, $false, "ሰላምታ በኋላ፣ ፊል ሶስት ልጆች አሉኝ አለ፣ ከዚያ", 2) | Out-Null
$d.Content.Find.Execute("Mike, astonished, asks: 'How old are they?' Fil,", $true, $false, $false, $false, $false, $true, 1, $false, "ሚኪ እጅግ ተገርሞ፣ ጠየቀው፡ ስንት አመታቸው ነወ? ፊል", 2) | Out-Null
$d.Content.Find.Execute("being a playful mathematician, answers", $true, $false, $false, $false, $false, $true, 1, $false, "እንደ ተጨዋች የሒሳብ ሊቅ፣ መለሰ", 2) | Out-Null
$d.Content.Find.Execute("'You tell me! I'll give you a hint: if you", $true, $false, $false, $false, $false, $true, 1, $false, "“ትነግርኛለህ! ፍንጭ እሰጥሃለሁ፡ አነተ", 2) | Out-Null
$d.Content.Find.Execute("multiply the three ages together you", $true, $false, $false, $false, $false, $true, 1, $false, "የሶስቱን ልጆች እድሜ አንድ ላይ ስታባዛው", 2) | Out-Null
$d.Content.Find.Execute("get 36.' Mike takes sometimes to think", $true, $false, $false, $false, $false, $true, 1, $false, "36 ታገኛለህ፡፡” ሚኪ ለማሰብ የተወሰነ ጊዜ ወሰደና", 2) | Out-Null
$d.Content.Find.Execute("and says: 'I'm sorry Fil, but I do need", $true, $false, $false, $false, $false, $true, 1, $false, "እንድህ አለው፡ አዝናለሁ ፊል፣ ነገር ግን", 2) | Out-Null
$d.Content.Find.Execute("another hint. So Fil tells Mike:", $true, $false, $false, $false, $false, $true, 1, $false, "ሌላ  ፍንጭ እፈልጋለሁ፡፡ ስለዚህ ፊል ለሚኪ ነገረው፡-", 2) | Out-Null
$d.Content.Find.Execute("'Yes, sure, here it is: if you had up to", $true, $false, $false, $false, $false, $true, 1, $false, "“አወ፣ በሚገባ፣ ይኸውልህ፡- ብታነሳለት እስከ", 2) | Out-Null
$d.Content.Find.Execute("three ages you get the number of math", $true, $false, $false, $false, $false, $true, 1, $false, "ሶስት አመት የምታገኘው ቁጥር በሒሳብ", 2) | Out-Null
$d.Content.Find.Execute("00:01:28,000 --> 00:01:31,000", $true, $false, $false, $false, $false, $true, 1, $false, "ሶስት አመት የምታገኘው ቁጥር በሒሳብ", 2) | Out-Null
$d.Content.Find.Execute("papers we publish together. Do you remember it?'", $true, $false, $false, $false, $false, $true, 1, $false, "በጋራ ያሳተምናቸውን  ህትመቶችን ነው፡፡ አስታወስከው?”", 2) | Out-Null
$d.Content.Find.Execute("'Yes I do remember How many, but still", $true, $false, $false, $false, $false, $true, 1, $false, "“አወ፣ ስንት እንደሆነ አስታወስኩት፣ ነገር ግን አሁንም", 2) | Out-Null
$d.Content.Find.Execute("I do not have enough information! I need", $true, $false, $false, $false, $false, $true, 1, $false, "በቂ መረጃ አላገኘሁም! እኔ", 2) | Out-Null
$d.Content.Find.Execute("at least one more.' Fil says: 'Yes don't", $true, $false, $false, $false, $false, $true, 1, $false, "አንድ ተጨማሪ እፈልጋለሁ፡፡” ፊል እንድህ አለ፡- “አወ", 2) | Out-Null
$d.Content.Find.Execute("worry but this is the last one:", $true, $false, $false, $false, $false, $true, 1, $false, "አትጨነቅ ግን ይህ የመጨረሻ ነው፡-", 2) | Out-Null
$d.Content.Find.Execute("The youngest one has blues eyes.' And", $true, $false, $false, $false, $false, $true, 1, $false, "ትንሹ ሴሚያዊ አይን አለው::”  እና ", 2) | Out-Null
$d.Content.Find.Execute("suddenly Mike gets the answer. You", $true, $false, $false, $false, $false, $true, 1, $false, "በቅጽበት ሚኪ መልሱን አገኘው፡፡ አንተ", 2) | Out-Null
$d.Content.Find.Execute("hear the conversation but you don't know", $true, $false, $false, $false, $false, $true, 1, $false, "ምልልሱን ሰምተሀል ግን አላወከውም", 2) | Out-Null
$d.Content.Find.Execute("how many papers they published together.", $true, $false, $false, $false, $false, $true, 1, $false, "ስንት ህትመት አብረው እንዳሳተሙ፡፡", 2) | Out-Null
$d.Content.Find.Execute("However, you do want to know the ages of", $true, $false, $false, $false, $false, $true, 1, $false, "ነገር ግን እድሜያቸውን ለማወቅ ፈልገሀል", 2) | Out-Null
$d.Content.Find.Execute("the three children. Can you figure them", $true, $false, $false, $false, $false, $true, 1, $false, "የሶስቱን ልጆች፡፡ አሁን በግልፅ", 2) | Out-Null
$d.Content.Find.Execute("out?", $true, $false, $false, $false, $false, $true, 1, $false, "ታዩህ?", 2) | Out-Null

Write-Output "Replacements applied"
